$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2280.9167
$ws.Range("I15").Value = 2280.9167
$ws.Range("K15").Value = 6842.750100000001
$ws.Range("M15").Value = -6673.750100000001
# Row 121
$ws.Range("H121").Value = 777.8261
$ws.Range("J121").Value = 802.381
$ws.Range("L121").Value = 2407.143
$ws.Range("N121").Value = -5901.143

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3453.6667
$ws.Range("I2").Value = 3844.4
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 3844.4
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -3731.4
$ws.Range("N2").Value = -1726
# Row 61
$ws.Range("H61").Value = 2035.2894
$ws.Range("I61").Value = 1861.0667
$ws.Range("J61").Value = 2688.625
$ws.Range("K61").Value = 1861.0667
$ws.Range("L61").Value = 2688.625
$ws.Range("M61").Value = -1649.0667
$ws.Range("N61").Value = -3112.625
# Row 110
$ws.Range("H110").Value = 1904.3334
$ws.Range("I110").Value = 1850
$ws.Range("J110").Value = 2013
$ws.Range("K110").Value = 1850
$ws.Range("L110").Value = 2013
$ws.Range("M110").Value = 195
$ws.Range("N110").Value = -6103
# Row 116
$ws.Range("H116").Value = 3453.6667
$ws.Range("I116").Value = 3844.4
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 3844.4
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = -1550.4
$ws.Range("N116").Value = -6088
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 136
$ws.Range("H136").Value = 2035.2894
$ws.Range("I136").Value = 1861.0667
$ws.Range("J136").Value = 2688.625
$ws.Range("K136").Value = 5583.2001
$ws.Range("L136").Value = 8065.875
$ws.Range("M136").Value = -3033.2001
$ws.Range("N136").Value = -13165.875

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3453.6667
$ws.Range("I3").Value = 3844.4
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 3844.4
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -3730.4
$ws.Range("N3").Value = -1728

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 3057.516
$ws.Range("I134").Value = 3032.8845
$ws.Range("J134").Value = 3185.6
$ws.Range("K134").Value = 9098.6535
$ws.Range("L134").Value = 9556.799999999999
$ws.Range("M134").Value = -6563.6535
$ws.Range("N134").Value = -14626.8

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 16666996
$ws.Range("I97").Value = 16666996
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 50000988
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -50000492
$ws.Range("N97").ClearContents()
# Row 101
$ws.Range("H101").Value = 10750.25
$ws.Range("J101").Value = 10750.25
$ws.Range("L101").Value = 32250.75
$ws.Range("N101").Value = -37118.75
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 112
$ws.Range("H112").Value = 3105.2632
$ws.Range("J112").Value = 3222.2222
$ws.Range("L112").Value = 9666.6666
$ws.Range("N112").Value = -11882.6666
# Row 121
$ws.Range("H121").Value = 902.5
$ws.Range("I121").Value = 276.125
$ws.Range("J121").Value = 1153.05
$ws.Range("K121").Value = 828.375
$ws.Range("L121").Value = 3459.15
$ws.Range("M121").Value = 481.625
$ws.Range("N121").Value = -6079.15
# Row 129
$ws.Range("H129").Value = 1638.138
$ws.Range("I129").Value = 1111
$ws.Range("J129").Value = 1915.579
$ws.Range("K129").Value = 3333
$ws.Range("L129").Value = 5746.737
$ws.Range("M129").Value = 1667
$ws.Range("N129").Value = -15746.737

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4312.0303
$ws.Range("J132").Value = 3707.0356
$ws.Range("L132").Value = 11121.1068
$ws.Range("N132").Value = -16181.1068

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6538038.5
$ws.Range("I22").Value = 18519358
$ws.Range("J22").Value = 2772.7273
$ws.Range("K22").Value = 18519358
$ws.Range("L22").Value = 2772.7273
$ws.Range("M22").Value = -18519063
$ws.Range("N22").Value = -3362.7273
# Row 27
$ws.Range("H27").Value = 6538038.5
$ws.Range("I27").Value = 18519358
$ws.Range("J27").Value = 2772.7273
$ws.Range("K27").Value = 18519358
$ws.Range("L27").Value = 2772.7273
$ws.Range("M27").Value = -18519251
$ws.Range("N27").Value = -2986.7273
# Row 68
$ws.Range("H68").Value = 37039130
$ws.Range("I68").Value = 1867.3158
$ws.Range("J68").Value = 125002616
$ws.Range("K68").Value = 1867.3158
$ws.Range("L68").Value = 125002616
$ws.Range("M68").Value = -1118.3158
$ws.Range("N68").Value = -125004114
# Row 71
$ws.Range("H71").Value = 37039130
$ws.Range("I71").Value = 1867.3158
$ws.Range("J71").Value = 125002616
$ws.Range("K71").Value = 9336.579
$ws.Range("L71").Value = 625013080
$ws.Range("M71").Value = -5592.579
$ws.Range("N71").Value = -625020568
# Row 82
$ws.Range("H82").Value = 8314.591
$ws.Range("I82").Value = 7843.2
$ws.Range("J82").Value = 9324.714
$ws.Range("K82").Value = 7843.2
$ws.Range("L82").Value = 9324.714
$ws.Range("M82").Value = -7482.2
$ws.Range("N82").Value = -10046.714
# Row 85
$ws.Range("H85").Value = 8314.591
$ws.Range("I85").Value = 7843.2
$ws.Range("J85").Value = 9324.714
$ws.Range("K85").Value = 7843.2
$ws.Range("L85").Value = 9324.714
$ws.Range("M85").Value = -6595.2
$ws.Range("N85").Value = -11820.714
# Row 100
$ws.Range("H100").Value = 1842.3
$ws.Range("I100").Value = 1467.1666
$ws.Range("J100").Value = 2405
$ws.Range("K100").Value = 1467.1666
$ws.Range("L100").Value = 2405
$ws.Range("M100").Value = -926.1666
$ws.Range("N100").Value = -3487
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1454.0714
$ws.Range("I122").Value = 1205.1818
$ws.Range("J122").Value = 2366.6667
$ws.Range("K122").Value = 3615.5454
$ws.Range("L122").Value = 7100.000100000001
$ws.Range("M122").Value = -1165.5454
$ws.Range("N122").Value = -12000.0001
# Row 132
$ws.Range("H132").Value = 3368.9375
$ws.Range("I132").Value = 2992.5
$ws.Range("J132").Value = 4498.25
$ws.Range("K132").Value = 8977.5
$ws.Range("L132").Value = 13494.75
$ws.Range("M132").Value = -6447.5
$ws.Range("N132").Value = -18554.75
